# RandomFloat66_HeapSortTimes.csv.xlsx — "Fix sorting and generate viable xlsx and charts"
#
# The two outlier Avg_Time_ms measurements for the smallest file sizes
# (5000 and 10000 rows) were corrected after fixing the sort routine.
# Update the underlying data cells; the scatter chart ("... Performance")
# plots Data!$D$2:$D$8 vs Data!$E$2:$E$8 so it reflects the corrected
# numbers the next time Excel recalculates/refreshes the chart cache.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.95313779
$ws.Range("D3").Value = 2.3136646

# Force a full recalculation so any formulas / chart caches relying on
# these cells are refreshed against the new values.
$wb.RefreshAll()
$excel.CalculateFullRebuild()
